{"js": "// Remove the \"Ver no Jupiter Salvar em pdf Salvar em docx\" paragraph, the\n// \"\u00a9 2020 . Contact: ...\" footer paragraph, and the blank paragraph that\n// separates them from the preceding \"LOB1019: F\u00edsica II (Requisito fraco)\"\n// line (build-site regeneration dropped these trailing Jekyll-site lines).\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nconst items = paragraphs.items;\n\n// Find the \"LOB1019: F\u00edsica II (Requisito fraco)\" paragraph; the three\n// paragraphs that immediately follow it (blank separator, \"Ver no Jupiter\n// ...\", and the \"\u00a9 2020 ...\" copyright line) are the ones being removed.\nlet anchorIndex = -1;\nfor (let i = 0; i < items.length; i++) {\n  if (items[i].text.indexOf(\"LOB1019\") !== -1) {\n    anchorIndex = i;\n    break;\n  }\n}\n\nif (anchorIndex !== -1) {\n  const toDelete = [];\n  for (let i = anchorIndex + 1; i < items.length && toDelete.length < 3; i++) {\n    toDelete.push(items[i]);\n  }\n  // Delete from the end backwards so earlier indices stay valid.\n  for (let i = toDelete.length - 1; i >= 0; i--) {\n    toDelete[i].delete();\n  }\n  await context.sync();\n}\n", "ps1": "# Remove the \"Ver no Jupiter Salvar em pdf Salvar em docx\" paragraph, the\n# \"\u00a9 2020 . Contact: ...\" footer paragraph, and the blank paragraph that\n# separates them from the preceding \"LOB1019: F\u00edsica II (Requisito fraco)\"\n# line (build-site regeneration dropped these trailing Jekyll-site lines).\n$d = $word.ActiveDocument\n\n$anchor = -1\n$count = $d.Paragraphs.Count\nfor ($i = 1; $i -le $count; $i++) {\n    if ($d.Paragraphs.Item($i).Range.Text -like \"*LOB1019*\") {\n        $anchor = $i\n        break\n    }\n}\n\nif ($anchor -ne -1) {\n    # The three paragraphs immediately after the anchor are the ones being\n    # removed. Delete starting from the furthest one so earlier indices stay\n    # valid while we work.\n    for ($k = 3; $k -ge 1; $k--) {\n        $idx = $anchor + $k\n        if ($idx -le $d.Paragraphs.Count) {\n            $d.Paragraphs.Item($idx).Range.Delete()\n        }\n    }\n}\n"}
